# Update "想去人数" (interest count) values in column F for the
# "展览" and "全部类型" sheets, reflecting a newer data snapshot.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 8425
    $ws.Range("F3").Value = 8017
    $ws.Range("F10").Value = 182
    $ws.Range("F12").Value = 723
    $ws.Range("F14").Value = 2077
    $ws.Range("F20").Value = 52
}
